{"js": "// Apply the dated worksheet refresh: update the header date and every\n// two-digit multiplication problem's operands according to the mapping\n// below. Each old value is unique in the document, so a direct search +\n// replace of the whole run text is safe and preserves formatting (the\n// run properties / fonts are untouched because we only replace the text).\nconst replacements = [\n  [\"2025-03-11 Tuesday\", \"2025-03-12 Wednesday\"],\n  [\"42\u00d748=\", \"84\u00d712=\"],\n  [\"22\u00d747=\", \"65\u00d723=\"],\n  [\"33\u00d748=\", \"12\u00d759=\"],\n  [\"32\u00d791=\", \"49\u00d781=\"],\n  [\"94\u00d724=\", \"48\u00d758=\"],\n  [\"45\u00d719=\", \"29\u00d743=\"],\n  [\"99\u00d765=\", \"34\u00d746=\"],\n  [\"44\u00d743=\", \"95\u00d753=\"],\n  [\"21\u00d787=\", \"53\u00d774=\"],\n  [\"12\u00d795=\", \"47\u00d762=\"],\n  [\"49\u00d759=\", \"33\u00d745=\"],\n  [\"67\u00d749=\", \"30\u00d752=\"],\n  [\"66\u00d714=\", \"91\u00d784=\"],\n  [\"50\u00d731=\", \"67\u00d766=\"],\n  [\"99\u00d784=\", \"18\u00d789=\"],\n  [\"26\u00d752=\", \"47\u00d758=\"],\n  [\"29\u00d766=\", \"87\u00d771=\"],\n  [\"49\u00d796=\", \"65\u00d764=\"],\n  [\"30\u00d718=\", \"22\u00d796=\"],\n  [\"22\u00d735=\", \"81\u00d775=\"],\n  [\"38\u00d793=\", \"25\u00d753=\"],\n  [\"81\u00d783=\", \"36\u00d763=\"],\n  [\"39\u00d738=\", \"25\u00d748=\"],\n  [\"59\u00d760=\", \"42\u00d737=\"],\n  [\"70\u00d740=\", \"16\u00d732=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the header date and every\n# two-digit multiplication problem's operands according to the mapping\n# below. Each old value is unique in the document, so Find/Replace over\n# the whole document content is safe and leaves run formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-11 Tuesday\", \"2025-03-12 Wednesday\"),\n    @(\"42\u00d748=\", \"84\u00d712=\"),\n    @(\"22\u00d747=\", \"65\u00d723=\"),\n    @(\"33\u00d748=\", \"12\u00d759=\"),\n    @(\"32\u00d791=\", \"49\u00d781=\"),\n    @(\"94\u00d724=\", \"48\u00d758=\"),\n    @(\"45\u00d719=\", \"29\u00d743=\"),\n    @(\"99\u00d765=\", \"34\u00d746=\"),\n    @(\"44\u00d743=\", \"95\u00d753=\"),\n    @(\"21\u00d787=\", \"53\u00d774=\"),\n    @(\"12\u00d795=\", \"47\u00d762=\"),\n    @(\"49\u00d759=\", \"33\u00d745=\"),\n    @(\"67\u00d749=\", \"30\u00d752=\"),\n    @(\"66\u00d714=\", \"91\u00d784=\"),\n    @(\"50\u00d731=\", \"67\u00d766=\"),\n    @(\"99\u00d784=\", \"18\u00d789=\"),\n    @(\"26\u00d752=\", \"47\u00d758=\"),\n    @(\"29\u00d766=\", \"87\u00d771=\"),\n    @(\"49\u00d796=\", \"65\u00d764=\"),\n    @(\"30\u00d718=\", \"22\u00d796=\"),\n    @(\"22\u00d735=\", \"81\u00d775=\"),\n    @(\"38\u00d793=\", \"25\u00d753=\"),\n    @(\"81\u00d783=\", \"36\u00d763=\"),\n    @(\"39\u00d738=\", \"25\u00d748=\"),\n    @(\"59\u00d760=\", \"42\u00d737=\"),\n    @(\"70\u00d740=\", \"16\u00d732=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
